$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price/volume table (and the two pairs of rows whose coins
# were re-ranked and swapped positions: rows 41/42 and rows 46/47).
# Numeric-looking text values (e.g. "1.003", "315.96") have their cell format
# forced to Text ("@") before the value is written, so Excel keeps them as the
# exact original text instead of silently converting them to numbers.

# Row 2: Bitcoin
$ws.Cells.Item(2, 4).Value = '28.273.65'
$ws.Cells.Item(2, 5).Value = '  +2.64%  '

# Row 3: Ethereum
$ws.Cells.Item(3, 4).Value = '1.878.59'
$ws.Cells.Item(3, 5).Value = '  +1.76%  '

# Row 4: TetherUSD
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.003'
$ws.Cells.Item(4, 5).Value = '  +0.19%  '

# Row 5: BNB
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '315.96'
$ws.Cells.Item(5, 5).Value = '  +0.86%  '

# Row 6: USDC
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '1.002'
$ws.Cells.Item(6, 5).Value = '  +0.07%  '

# Row 7: XRP
$ws.Cells.Item(7, 5).Value = '  +1.52%  '

# Row 8: Cardano
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.3731'
$ws.Cells.Item(8, 5).Value = '  +2.68%  '

# Row 9: Dogecoin
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.07419'
$ws.Cells.Item(9, 5).Value = '  +1.77%  '

# Row 10: Polygon
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.8850'
$ws.Cells.Item(10, 5).Value = '  +1.30%  '

# Row 11: Solana
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '21.15'
$ws.Cells.Item(11, 5).Value = '  +2.14%  '

# Row 12: WrappedEther
$ws.Cells.Item(12, 4).Value = '1.918.99'
$ws.Cells.Item(12, 5).Value = '  +3.18%  '

# Row 13: Polkadot
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '5.500'
$ws.Cells.Item(13, 5).Value = '  +3.13%  '

# Row 14: Chainlink
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '6.638'
$ws.Cells.Item(14, 5).Value = '  +1.98%  '

# Row 15: TRON
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.06989'
$ws.Cells.Item(15, 5).Value = '  +1.16%  '

# Row 16: BinanceUSD
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '1.005'
$ws.Cells.Item(16, 5).Value = '  +0.14%  '

# Row 17: Litecoin
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '81.35'
$ws.Cells.Item(17, 5).Value = '  +2.83%  '

# Row 18: ShibaInu
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.000009140'
$ws.Cells.Item(18, 5).Value = '  +3.12%  '

# Row 19: Dai
$ws.Cells.Item(19, 5).Value = '  -0.05%  '

# Row 20: Avalanche
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '15.65'
$ws.Cells.Item(20, 5).Value = '  +1.88%  '

# Row 21: WrappedBTC
$ws.Cells.Item(21, 4).Value = '28.576.16'
$ws.Cells.Item(21, 5).Value = '  +3.67%  '

# Row 22: Uniswap
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '5.085'
$ws.Cells.Item(22, 5).Value = '  +1.98%  '

# Row 23: Cosmos
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '10.96'
$ws.Cells.Item(23, 5).Value = '  +5.37%  '

# Row 24: WrappedliquidstakedEther2.0
$ws.Cells.Item(24, 4).Value = '2.225.17'
$ws.Cells.Item(24, 5).Value = '  +6.93%  '

# Row 25: Toncoin
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '1.972'
$ws.Cells.Item(25, 5).Value = '  -0.75%  '

# Row 26: Monero
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '154.39'
$ws.Cells.Item(26, 5).Value = '  +1.14%  '

# Row 27: EthereumClassic
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '18.83'
$ws.Cells.Item(27, 5).Value = '  -0.44%  '

# Row 28: InternetComputer(DFINITY)
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '5.410'
$ws.Cells.Item(28, 5).Value = '  +3.21%  '

# Row 29: BitcoinCash
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '117.51'
$ws.Cells.Item(29, 5).Value = '  -3.08%  '

# Row 30: LidoDAOToken
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '1.869'
$ws.Cells.Item(30, 5).Value = '  -0.62%  '

# Row 31: Stellar
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.08999'
$ws.Cells.Item(31, 5).Value = '  +1.43%  '

# Row 32: ImmutableX
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.7928'
$ws.Cells.Item(32, 5).Value = '  +3.92%  '

# Row 33: Filecoin
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '4.695'
$ws.Cells.Item(33, 5).Value = '  +3.06%  '

# Row 34: ARBITRUM
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.183'
$ws.Cells.Item(34, 5).Value = '  +7.67%  '

# Row 35: HuobiToken
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '2.950'
$ws.Cells.Item(35, 5).Value = '  +0.11%  '

# Row 36: Frax
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '1.002'
$ws.Cells.Item(36, 5).Value = '  +0.10%  '

# Row 37: TrustWalletToken
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '1.130'
$ws.Cells.Item(37, 5).Value = '  +3.66%  '

# Row 38: Hedera
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.05461'
$ws.Cells.Item(38, 5).Value = '  +2.22%  '

# Row 39: VeChain
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.01965'
$ws.Cells.Item(39, 5).Value = '  +1.77%  '

# Row 40: MXToken
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '2.892'
$ws.Cells.Item(40, 5).Value = '  +2.93%  '

# Row 41: TheSandbox
$ws.Cells.Item(41, 2).Value = 'TheSandbox'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.5175'
$ws.Cells.Item(41, 5).Value = '  +1.60%  '

# Row 42: Algorand
$ws.Cells.Item(42, 2).Value = 'Algorand'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.1689'
$ws.Cells.Item(42, 5).Value = '  +2.69%  '

# Row 43: FraxShare
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '6.908'
$ws.Cells.Item(43, 5).Value = '  +0.79%  '

# Row 44: Aptos
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '8.682'
$ws.Cells.Item(44, 5).Value = '  +5.07%  '

# Row 45: EnergySwap
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '10.62'
$ws.Cells.Item(45, 5).Value = '  +2.61%  '

# Row 46: Cronos
$ws.Cells.Item(46, 2).Value = 'Cronos'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.06606'
$ws.Cells.Item(46, 5).Value = '  +1.19%  '

# Row 47: Decentraland
$ws.Cells.Item(47, 2).Value = 'Decentraland'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.4765'
$ws.Cells.Item(47, 5).Value = '  +0.56%  '

# Row 48: Quant
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '106.41'
$ws.Cells.Item(48, 5).Value = '  +1.82%  '

# Row 49: PaxDollar
$ws.Cells.Item(49, 5).Value = '  +0.07%  '

# Row 50: NEARProtocol
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.658'
$ws.Cells.Item(50, 5).Value = '  +2.28%  '

# Row 51: RenderToken
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '1.830'
$ws.Cells.Item(51, 5).Value = '  +6.00%  '

Write-Host "Updated cryptos list on Wed May 10 16:49:43 UTC 2023 with GitHub Actions"
